$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.215.30"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "1.860.60"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.73%  "

$ws.Range("D5").Value = "'239.53"
$ws.Range("E5").Value = "  +3.68%  "

$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D8").Value = "'42.31"
$ws.Range("E8").Value = "  +6.29%  "

$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").Value = "'0.0694"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").Value = "'0.0990"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "2.130.07"

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.868.31"
$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.48"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("D17").Value = "35.227.87"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").Value = "'69.93"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "'241.44"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'2.26"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'169.65"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").Value = "'1.89"
$ws.Range("E26").Value = "  +24.54%  "

$ws.Range("D27").Value = "'8.08"
$ws.Range("E27").Value = "  +3.72%  "

$ws.Range("D28").Value = "'17.70"
$ws.Range("E28").Value = "  +1.82%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D30").Value = "'0.0562"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").Value = "'4.02"
$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("E33").Value = "  +28.99%  "

$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").Value = "'2.05"
$ws.Range("E35").Value = "  +9.59%  "

$ws.Range("D36").Value = "'0.817"
$ws.Range("E36").Value = "  +17.77%  "

$ws.Range("E37").Value = "  +7.51%  "

$ws.Range("E38").Value = "  +3.58%  "

$ws.Range("D39").Value = "'0.0203"
$ws.Range("E39").Value = "  +4.77%  "

$ws.Range("D40").Value = "'90.19"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").Value = "1.345.75"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D42").Value = "'0.0601"
$ws.Range("E42").Value = "  +15.02%  "

$ws.Range("D43").Value = "'15.14"
$ws.Range("E43").Value = "  +3.08%  "

$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'12.49"
$ws.Range("E46").Value = "  +46.93%  "

$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "'6.59"
$ws.Range("E48").Value = "  +5.13%  "

$ws.Range("D49").Value = "2.045.58"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("E50").Value = "  +3.28%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'3.43"
$ws.Range("E51").Value = "  +4.17%  "
